$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / title changes: October -> November ---
$ws.Range("A2").Value = "November 2016 and November 2015 (Thousand Tons)"

# Row 5 header cells hold month/year text (e.g. "October 2016") which Excel
# auto-parses as a date if assigned directly. Force text format, assign, then
# restore the original number format so styling is unaffected.
$headerCells2016 = @("B5","E5","G5","I5","K5")
foreach ($addr in $headerCells2016) {
  $rng = $ws.Range($addr)
  $fmt = $rng.NumberFormat
  $rng.NumberFormat = "@"
  $rng.Value = "November 2016"
  $rng.NumberFormat = $fmt
}
$headerCells2015 = @("C5","F5","H5","J5","L5")
foreach ($addr in $headerCells2015) {
  $rng = $ws.Range($addr)
  $fmt = $rng.NumberFormat
  $rng.NumberFormat = "@"
  $rng.Value = "November 2015"
  $rng.NumberFormat = $fmt
}

# --- Data cell updates ---
$ws.Range("B6").Value = 287
$ws.Range("C6").Value = 316
$ws.Range("D6").Value = -0.092
$ws.Range("G6").Value = 271
$ws.Range("H6").Value = 301
$ws.Range("I6").Value = 16
$ws.Range("J6").Value = 16
$ws.Range("B7").Value = 104
$ws.Range("C7").Value = 112
$ws.Range("D7").Value = -0.069
$ws.Range("G7").Value = 104
$ws.Range("H7").Value = 112
$ws.Range("B8").Value = 23
$ws.Range("C8").Value = 26
$ws.Range("D8").Value = -0.12
$ws.Range("H8").Value = 10
$ws.Range("I8").Value = 16
$ws.Range("J8").Value = 16
$ws.Range("B9").Value = 152
$ws.Range("C9").Value = 168
$ws.Range("D9").Value = -0.097
$ws.Range("G9").Value = 152
$ws.Range("H9").Value = 168
$ws.Range("B10").Value = 9
$ws.Range("D10").Value = -0.18
$ws.Range("G10").Value = 9
$ws.Range("B13").Value = 413
$ws.Range("C13").Value = 450
$ws.Range("D13").Value = -0.083
$ws.Range("G13").Value = 330
$ws.Range("H13").Value = 355
$ws.Range("I13").Value = 83
$ws.Range("J13").Value = 95
$ws.Range("B14").Value = 107
$ws.Range("C14").Value = 117
$ws.Range("D14").Value = -0.083
$ws.Range("G14").Value = 82
$ws.Range("H14").Value = 90
$ws.Range("I14").Value = 26
$ws.Range("J14").Value = 28
$ws.Range("B15").Value = 153
$ws.Range("C15").Value = 170
$ws.Range("D15").Value = -0.099
$ws.Range("G15").Value = 118
$ws.Range("H15").Value = 127
$ws.Range("I15").Value = 36
$ws.Range("J15").Value = 43
$ws.Range("B16").Value = 152
$ws.Range("C16").Value = 163
$ws.Range("D16").Value = -0.066
$ws.Range("G16").Value = 130
$ws.Range("H16").Value = 139
$ws.Range("I16").Value = 22
$ws.Range("J16").Value = 24
$ws.Range("D17").Value = -0.14
$ws.Range("I17").Value = 16
$ws.Range("J17").Value = 20
$ws.Range("D19").Value = -0.21
$ws.Range("D20").Value = -0.16
$ws.Range("D22").Value = 0.03
$ws.Range("B23").Value = 52
$ws.Range("C23").Value = 54
$ws.Range("D23").Value = -0.031
$ws.Range("E23").Value = 36
$ws.Range("F23").Value = 31
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = 21
$ws.Range("I23").Value = 1
$ws.Range("B26").Value = 52
$ws.Range("C26").Value = 54
$ws.Range("D26").Value = -0.031
$ws.Range("E26").Value = 36
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 15
$ws.Range("H26").Value = 21
$ws.Range("I26").Value = 1
$ws.Range("B31").Value = 426
$ws.Range("C31").Value = 444
$ws.Range("D31").Value = -0.04
$ws.Range("G31").Value = 394
$ws.Range("H31").Value = 410
$ws.Range("I31").Value = 31
$ws.Range("B34").Value = 274
$ws.Range("C34").Value = 284
$ws.Range("D34").Value = -0.036
$ws.Range("G34").Value = 274
$ws.Range("H34").Value = 284
$ws.Range("B36").Value = 67
$ws.Range("D36").Value = -0.028
$ws.Range("G36").Value = 67
$ws.Range("B39").Value = 85
$ws.Range("C39").Value = 91
$ws.Range("D39").Value = -0.064
$ws.Range("G39").Value = 54
$ws.Range("H39").Value = 57
$ws.Range("I39").Value = 31
$ws.Range("B46").Value = 0.43
$ws.Range("D46").Value = -0.52
$ws.Range("K46").Value = 0.43
$ws.Range("B49").Value = 0.43
$ws.Range("D49").Value = -0.52
$ws.Range("K49").Value = 0.43
$ws.Range("B51").Value = 0.13
$ws.Range("D51").Value = -0.27
$ws.Range("G51").Value = 0.13
$ws.Range("B58").Value = 0.13
$ws.Range("D58").Value = -0.27
$ws.Range("G58").Value = 0.13
$ws.Range("B60").Value = 59
$ws.Range("C60").Value = 66
$ws.Range("D60").Value = -0.11
$ws.Range("G60").Value = 59
$ws.Range("H60").Value = 66
$ws.Range("B61").Value = 41
$ws.Range("C61").Value = 43
$ws.Range("D61").Value = -0.064
$ws.Range("G61").Value = 41
$ws.Range("H61").Value = 43
$ws.Range("C62").Value = 11
$ws.Range("H62").Value = 11
$ws.Range("B63").Value = 11
$ws.Range("C63").Value = 12
$ws.Range("D63").Value = -0.071
$ws.Range("G63").Value = 11
$ws.Range("H63").Value = 12
$ws.Range("B64").Value = 36
$ws.Range("C64").Value = 35
$ws.Range("D64").Value = 0.047
$ws.Range("I64").Value = 36
$ws.Range("J64").Value = 35
$ws.Range("B66").Value = 36
$ws.Range("C66").Value = 35
$ws.Range("D66").Value = 0.047
$ws.Range("I66").Value = 36
$ws.Range("J66").Value = 35
$ws.Range("B67").Value = 1294
$ws.Range("C67").Value = 1389
$ws.Range("D67").Value = -0.069
$ws.Range("E67").Value = 39
$ws.Range("F67").Value = 34
$ws.Range("G67").Value = 1069
$ws.Range("H67").Value = 1153
$ws.Range("I67").Value = 184
$ws.Range("J67").Value = 202
$ws.Range("K67").Value = 0.43
